$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 451.375
$ws.Cells.Item(2, 9).Value = 290.25
$ws.Cells.Item(2, 10).Value = 612.5
$ws.Cells.Item(2, 11).Value = 290.25
$ws.Cells.Item(2, 12).Value = 612.5
$ws.Cells.Item(2, 13).Value = -177.25
$ws.Cells.Item(2, 14).Value = -838.5
$ws.Cells.Item(10, 8).Value = 0
$ws.Cells.Item(10, 10).Value = 0
$ws.Cells.Item(10, 12).Value = 0
$ws.Cells.Item(10, 14).Value = $null
$ws.Cells.Item(15, 8).Value = 1280.4231
$ws.Cells.Item(15, 9).Value = 1280.4231
$ws.Cells.Item(15, 11).Value = 3841.2693
$ws.Cells.Item(15, 13).Value = -3672.2693
$ws.Cells.Item(19, 8).Value = 1169.6471
$ws.Cells.Item(19, 9).Value = 437
$ws.Cells.Item(19, 10).Value = 1395.0769
$ws.Cells.Item(19, 11).Value = 437
$ws.Cells.Item(19, 12).Value = 1395.0769
$ws.Cells.Item(19, 13).Value = -262
$ws.Cells.Item(19, 14).Value = -1745.0769
$ws.Cells.Item(43, 8).Value = 2567291.5
$ws.Cells.Item(43, 10).Value = 4161.6665
$ws.Cells.Item(43, 12).Value = 4161.6665
$ws.Cells.Item(43, 14).Value = -4299.6665
$ws.Cells.Item(51, 8).Value = 7073.857
$ws.Cells.Item(51, 10).Value = 8007.2144
$ws.Cells.Item(51, 12).Value = 8007.2144
$ws.Cells.Item(51, 14).Value = -8975.214400000001
$ws.Cells.Item(64, 8).Value = 60390356
$ws.Cells.Item(64, 10).Value = 4998
$ws.Cells.Item(64, 12).Value = 4998
$ws.Cells.Item(64, 14).Value = -5494
$ws.Cells.Item(67, 8).Value = 60390356
$ws.Cells.Item(67, 10).Value = 4998
$ws.Cells.Item(67, 12).Value = 4998
$ws.Cells.Item(67, 14).Value = -6714
$ws.Cells.Item(70, 8).Value = 5891.769
$ws.Cells.Item(70, 10).Value = 7100
$ws.Cells.Item(70, 12).Value = 21300
$ws.Cells.Item(70, 14).Value = -21840
$ws.Cells.Item(73, 8).Value = 5891.769
$ws.Cells.Item(73, 10).Value = 7100
$ws.Cells.Item(73, 12).Value = 21300
$ws.Cells.Item(73, 14).Value = -23172
$ws.Cells.Item(101, 8).Value = 1036
$ws.Cells.Item(101, 9).Value = 1024
$ws.Cells.Item(101, 11).Value = 3072
$ws.Cells.Item(101, 13).Value = -1450
$ws.Cells.Item(104, 8).Value = 779.75
$ws.Cells.Item(104, 9).Value = 773
$ws.Cells.Item(104, 11).Value = 2319
$ws.Cells.Item(104, 13).Value = -572
$ws.Cells.Item(111, 8).Value = 4706.1055
$ws.Cells.Item(111, 9).Value = 4130.143
$ws.Cells.Item(111, 10).Value = 6318.8
$ws.Cells.Item(111, 11).Value = 12390.429
$ws.Cells.Item(111, 12).Value = 18956.4
$ws.Cells.Item(111, 13).Value = -9323.429
$ws.Cells.Item(111, 14).Value = -25090.4
$ws.Cells.Item(112, 8).Value = 3044.82
$ws.Cells.Item(112, 10).Value = 3128.9148
$ws.Cells.Item(112, 12).Value = 9386.7444
$ws.Cells.Item(112, 14).Value = -11602.7444
$ws.Cells.Item(121, 8).Value = 3977.524
$ws.Cells.Item(121, 10).Value = 3977.524
$ws.Cells.Item(121, 12).Value = 11932.572
$ws.Cells.Item(121, 14).Value = -15426.572
$ws.Cells.Item(129, 8).Value = 1764.7778
$ws.Cells.Item(129, 9).Value = 1017.4
$ws.Cells.Item(129, 10).Value = 2699
$ws.Cells.Item(129, 11).Value = 3052.2
$ws.Cells.Item(129, 12).Value = 8097
$ws.Cells.Item(129, 13).Value = 1947.8
$ws.Cells.Item(129, 14).Value = -18097
$ws.Cells.Item(135, 8).Value = 1862.7561
$ws.Cells.Item(135, 9).Value = 556.6111
$ws.Cells.Item(135, 11).Value = 5009.4999
$ws.Cells.Item(135, 13).Value = -2474.4999
$ws.Cells.Item(137, 8).Value = 14494951
$ws.Cells.Item(137, 9).Value = 1875.1333
$ws.Cells.Item(137, 10).Value = 41669468
$ws.Cells.Item(137, 11).Value = 5625.3999
$ws.Cells.Item(137, 12).Value = 125008404
$ws.Cells.Item(137, 13).Value = -3075.3999
$ws.Cells.Item(137, 14).Value = -125013504
$ws.Cells.Item(138, 8).Value = 5171.6665
$ws.Cells.Item(138, 9).Value = 1882.6786
$ws.Cells.Item(138, 10).Value = 6657.016
$ws.Cells.Item(138, 11).Value = 5648.0358
$ws.Cells.Item(138, 12).Value = 19971.048
$ws.Cells.Item(138, 13).Value = -508.0357999999997
$ws.Cells.Item(138, 14).Value = -30251.048
$ws.Cells.Item(141, 8).Value = 2956.4062
$ws.Cells.Item(141, 9).Value = 3324.7036
$ws.Cells.Item(141, 10).Value = 967.6
$ws.Cells.Item(141, 11).Value = 9974.110799999999
$ws.Cells.Item(141, 12).Value = 2902.8
$ws.Cells.Item(141, 13).Value = -4794.110799999999
$ws.Cells.Item(141, 14).Value = -13262.8

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(17, 8).Value = 9975
$ws.Cells.Item(17, 10).Value = 9975
$ws.Cells.Item(17, 12).Value = 9975
$ws.Cells.Item(17, 14).Value = -10321
$ws.Cells.Item(32, 8).Value = 16087.189
$ws.Cells.Item(32, 9).Value = 15971.64
$ws.Cells.Item(32, 10).Value = 16809.375
$ws.Cells.Item(32, 11).Value = 15971.64
$ws.Cells.Item(32, 12).Value = 16809.375
$ws.Cells.Item(32, 13).Value = -15684.64
$ws.Cells.Item(32, 14).Value = -17383.375
$ws.Cells.Item(45, 8).Value = 3349.4614
$ws.Cells.Item(45, 9).Value = 2598.4285
$ws.Cells.Item(45, 10).Value = 4225.6665
$ws.Cells.Item(45, 11).Value = 2598.4285
$ws.Cells.Item(45, 12).Value = 4225.6665
$ws.Cells.Item(45, 13).Value = -2221.4285
$ws.Cells.Item(45, 14).Value = -4979.6665
$ws.Cells.Item(74, 8).Value = 17859068
$ws.Cells.Item(74, 9).Value = 25001320
$ws.Cells.Item(74, 10).Value = 3437
$ws.Cells.Item(74, 11).Value = 25001320
$ws.Cells.Item(74, 12).Value = 3437
$ws.Cells.Item(74, 13).Value = -25000446
$ws.Cells.Item(74, 14).Value = -5185
$ws.Cells.Item(77, 8).Value = 17859068
$ws.Cells.Item(77, 9).Value = 25001320
$ws.Cells.Item(77, 10).Value = 3437
$ws.Cells.Item(77, 11).Value = 125006600
$ws.Cells.Item(77, 12).Value = 17185
$ws.Cells.Item(77, 13).Value = -125002232
$ws.Cells.Item(77, 14).Value = -25921
$ws.Cells.Item(110, 8).Value = 39308136
$ws.Cells.Item(110, 9).Value = 2917266
$ws.Cells.Item(110, 11).Value = 2917266
$ws.Cells.Item(110, 13).Value = -2915221
$ws.Cells.Item(112, 8).Value = 67999.39999999999
$ws.Cells.Item(112, 10).Value = 67999.39999999999
$ws.Cells.Item(112, 12).Value = 67999.39999999999
$ws.Cells.Item(112, 14).Value = -70953.39999999999
$ws.Cells.Item(122, 8).Value = 4398.593
$ws.Cells.Item(122, 9).Value = 3148.318
$ws.Cells.Item(122, 11).Value = 9444.954000000002
$ws.Cells.Item(122, 13).Value = -6994.954000000002
$ws.Cells.Item(132, 8).Value = 25381.45
$ws.Cells.Item(132, 9).Value = 31940.105
$ws.Cells.Item(132, 11).Value = 95820.315
$ws.Cells.Item(132, 13).Value = -93290.315
$ws.Cells.Item(133, 8).Value = 82994.5
$ws.Cells.Item(133, 10).Value = 82994.5
$ws.Cells.Item(133, 12).Value = 82994.5
$ws.Cells.Item(133, 14).Value = -88054.5
$ws.Cells.Item(139, 8).Value = 72461
$ws.Cells.Item(139, 10).Value = 72461
$ws.Cells.Item(139, 12).Value = 72461
$ws.Cells.Item(139, 14).Value = -82741

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(22, 8).Value = 821.2857
$ws.Cells.Item(22, 9).Value = 958.75
$ws.Cells.Item(22, 10).Value = 638
$ws.Cells.Item(22, 11).Value = 958.75
$ws.Cells.Item(22, 12).Value = 638
$ws.Cells.Item(22, 13).Value = -785.75
$ws.Cells.Item(22, 14).Value = -984
$ws.Cells.Item(86, 8).Value = 2564.3333
$ws.Cells.Item(86, 9).Value = 2286.111
$ws.Cells.Item(86, 11).Value = 2286.111
$ws.Cells.Item(86, 13).Value = -1163.111
$ws.Cells.Item(89, 8).Value = 2564.3333
$ws.Cells.Item(89, 9).Value = 2286.111
$ws.Cells.Item(89, 11).Value = 11430.555
$ws.Cells.Item(89, 13).Value = -5814.555
$ws.Cells.Item(105, 8).Value = 46876664
$ws.Cells.Item(105, 9).Value = 46876664
$ws.Cells.Item(105, 10).Value = 0
$ws.Cells.Item(105, 11).Value = 46876664
$ws.Cells.Item(105, 12).Value = 0
$ws.Cells.Item(105, 13).Value = -46874917
$ws.Cells.Item(105, 14).Value = $null
$ws.Cells.Item(107, 8).Value = 1699.7
$ws.Cells.Item(107, 9).Value = 1677.4445
$ws.Cells.Item(107, 11).Value = 1677.4445
$ws.Cells.Item(107, 13).Value = 242.5554999999999
$ws.Cells.Item(134, 8).Value = 2457.8572
$ws.Cells.Item(134, 9).Value = 2457.8572
$ws.Cells.Item(134, 11).Value = 7373.571599999999
$ws.Cells.Item(134, 13).Value = -4838.571599999999

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 2899.5
$ws.Cells.Item(16, 9).Value = 2899.5
$ws.Cells.Item(16, 10).Value = 0
$ws.Cells.Item(16, 11).Value = 2899.5
$ws.Cells.Item(16, 12).Value = 0
$ws.Cells.Item(16, 13).Value = -2612.5
$ws.Cells.Item(16, 14).Value = $null
$ws.Cells.Item(31, 8).Value = 52636348
$ws.Cells.Item(31, 9).Value = 76925510
$ws.Cells.Item(31, 10).Value = 9833.333000000001
$ws.Cells.Item(31, 11).Value = 76925510
$ws.Cells.Item(31, 12).Value = 9833.333000000001
$ws.Cells.Item(31, 13).Value = -76925215
$ws.Cells.Item(31, 14).Value = -10423.333
$ws.Cells.Item(34, 8).Value = 52636348
$ws.Cells.Item(34, 9).Value = 76925510
$ws.Cells.Item(34, 10).Value = 9833.333000000001
$ws.Cells.Item(34, 11).Value = 76925510
$ws.Cells.Item(34, 12).Value = 9833.333000000001
$ws.Cells.Item(34, 13).Value = -76925308
$ws.Cells.Item(34, 14).Value = -10237.333
$ws.Cells.Item(94, 8).Value = 2612.7646
$ws.Cells.Item(94, 10).Value = 3324.0833
$ws.Cells.Item(94, 12).Value = 3324.0833
$ws.Cells.Item(94, 14).Value = -4226.0833
$ws.Cells.Item(105, 8).Value = 2067887.4
$ws.Cells.Item(105, 9).Value = 3248204.2
$ws.Cells.Item(105, 10).Value = 2333
$ws.Cells.Item(105, 11).Value = 3248204.2
$ws.Cells.Item(105, 12).Value = 2333
$ws.Cells.Item(105, 13).Value = -3246457.2
$ws.Cells.Item(105, 14).Value = -5827
$ws.Cells.Item(113, 8).Value = 2899.5
$ws.Cells.Item(113, 9).Value = 2899.5
$ws.Cells.Item(113, 10).Value = 0
$ws.Cells.Item(113, 11).Value = 2899.5
$ws.Cells.Item(113, 12).Value = 0
$ws.Cells.Item(113, 13).Value = -729.5
$ws.Cells.Item(113, 14).Value = $null
$ws.Cells.Item(122, 8).Value = 4482.5293
$ws.Cells.Item(122, 9).Value = 2093.2856
$ws.Cells.Item(122, 11).Value = 6279.8568
$ws.Cells.Item(122, 13).Value = -3829.8568
$ws.Cells.Item(132, 8).Value = 34201440
$ws.Cells.Item(132, 9).Value = 43020400
$ws.Cells.Item(132, 10).Value = 27974.375
$ws.Cells.Item(132, 11).Value = 129061200
$ws.Cells.Item(132, 12).Value = 83923.125
$ws.Cells.Item(132, 13).Value = -129058670
$ws.Cells.Item(132, 14).Value = -88983.125
$ws.Cells.Item(134, 8).Value = 1718.5128
$ws.Cells.Item(134, 9).Value = 1730.3243
$ws.Cells.Item(134, 11).Value = 5190.9729
$ws.Cells.Item(134, 13).Value = -2655.9729
$ws.Cells.Item(135, 8).Value = 0
$ws.Cells.Item(135, 10).Value = 0
$ws.Cells.Item(135, 12).Value = 0
$ws.Cells.Item(135, 14).Value = $null

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(14, 8).Value = 289.18182
$ws.Cells.Item(14, 9).Value = 289.18182
$ws.Cells.Item(14, 11).Value = 867.54546
$ws.Cells.Item(14, 13).Value = -694.54546
$ws.Cells.Item(18, 8).Value = 7250
$ws.Cells.Item(18, 9).Value = 7000
$ws.Cells.Item(18, 10).Value = 7500
$ws.Cells.Item(18, 11).Value = 21000
$ws.Cells.Item(18, 12).Value = 22500
$ws.Cells.Item(18, 13).Value = -20831
$ws.Cells.Item(18, 14).Value = -22838
$ws.Cells.Item(22, 8).Value = 1550
$ws.Cells.Item(22, 9).Value = 0
$ws.Cells.Item(22, 10).Value = 1550
$ws.Cells.Item(22, 11).Value = 0
$ws.Cells.Item(22, 12).Value = 4650
$ws.Cells.Item(22, 13).Value = $null
$ws.Cells.Item(22, 14).Value = -4988
$ws.Cells.Item(27, 8).Value = 1550
$ws.Cells.Item(27, 9).Value = 0
$ws.Cells.Item(27, 10).Value = 1550
$ws.Cells.Item(27, 11).Value = 0
$ws.Cells.Item(27, 12).Value = 4650
$ws.Cells.Item(27, 13).Value = $null
$ws.Cells.Item(27, 14).Value = -4854
$ws.Cells.Item(34, 8).Value = 2953
$ws.Cells.Item(34, 10).Value = 1833.3334
$ws.Cells.Item(34, 12).Value = 5500.0002
$ws.Cells.Item(34, 14).Value = -5668.0002
$ws.Cells.Item(97, 8).Value = 409.66666
$ws.Cells.Item(97, 9).Value = 401.6
$ws.Cells.Item(97, 10).Value = 450
$ws.Cells.Item(97, 11).Value = 1204.8
$ws.Cells.Item(97, 12).Value = 1350
$ws.Cells.Item(97, 13).Value = -708.8000000000002
$ws.Cells.Item(97, 14).Value = -2342
$ws.Cells.Item(99, 8).Value = 5624.6665
$ws.Cells.Item(99, 10).Value = 6399.6
$ws.Cells.Item(99, 12).Value = 19198.8
$ws.Cells.Item(99, 14).Value = -23690.8
$ws.Cells.Item(119, 8).Value = 9209
$ws.Cells.Item(119, 9).Value = 777.8
$ws.Cells.Item(119, 10).Value = 19748
$ws.Cells.Item(119, 11).Value = 2333.4
$ws.Cells.Item(119, 12).Value = 59244
$ws.Cells.Item(119, 13).Value = 2504.6
$ws.Cells.Item(119, 14).Value = -68920
$ws.Cells.Item(124, 8).Value = 10805
$ws.Cells.Item(124, 9).Value = 1341.6666
$ws.Cells.Item(124, 10).Value = 25000
$ws.Cells.Item(124, 11).Value = 4024.9998
$ws.Cells.Item(124, 12).Value = 75000
$ws.Cells.Item(124, 13).Value = 885.0001999999999
$ws.Cells.Item(124, 14).Value = -84820
$ws.Cells.Item(126, 8).Value = 10566
$ws.Cells.Item(126, 9).Value = 943.3333
$ws.Cells.Item(126, 10).Value = 25000
$ws.Cells.Item(126, 11).Value = 2829.9999
$ws.Cells.Item(126, 12).Value = 75000
$ws.Cells.Item(126, 13).Value = 2110.0001
$ws.Cells.Item(126, 14).Value = -84880
$ws.Cells.Item(129, 8).Value = 1579.2354
$ws.Cells.Item(129, 10).Value = 2019.8
$ws.Cells.Item(129, 12).Value = 6059.4
$ws.Cells.Item(129, 14).Value = -16059.4
$ws.Cells.Item(131, 8).Value = 16506647
$ws.Cells.Item(131, 9).Value = 16667570
$ws.Cells.Item(131, 10).Value = 16471664
$ws.Cells.Item(131, 11).Value = 50002710
$ws.Cells.Item(131, 12).Value = 49414992
$ws.Cells.Item(131, 13).Value = -49997670
$ws.Cells.Item(131, 14).Value = -49425072
$ws.Cells.Item(137, 8).Value = 34002340
$ws.Cells.Item(137, 10).Value = 6252973
$ws.Cells.Item(137, 12).Value = 18758919
$ws.Cells.Item(137, 14).Value = -18769119

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 328.34482
$ws.Cells.Item(2, 9).Value = 136.66667
$ws.Cells.Item(2, 10).Value = 463.64706
$ws.Cells.Item(2, 11).Value = 136.66667
$ws.Cells.Item(2, 12).Value = 463.64706
$ws.Cells.Item(2, 13).Value = -23.66667000000001
$ws.Cells.Item(2, 14).Value = -689.64706
$ws.Cells.Item(12, 8).Value = 0
$ws.Cells.Item(12, 10).Value = 0
$ws.Cells.Item(12, 12).Value = 0
$ws.Cells.Item(12, 14).Value = $null
$ws.Cells.Item(80, 8).Value = 1863777.2
$ws.Cells.Item(80, 9).Value = 2774832.8
$ws.Cells.Item(80, 11).Value = 2774832.8
$ws.Cells.Item(80, 13).Value = -2773834.8
$ws.Cells.Item(83, 8).Value = 1863777.2
$ws.Cells.Item(83, 9).Value = 2774832.8
$ws.Cells.Item(83, 11).Value = 13874164
$ws.Cells.Item(83, 13).Value = -13869172
$ws.Cells.Item(99, 8).Value = 6707.8335
$ws.Cells.Item(99, 9).Value = 2049.4
$ws.Cells.Item(99, 11).Value = 2049.4
$ws.Cells.Item(99, 13).Value = 196.5999999999999
$ws.Cells.Item(113, 8).Value = 2588.875
$ws.Cells.Item(113, 9).Value = 2490.75
$ws.Cells.Item(113, 10).Value = 2687
$ws.Cells.Item(113, 11).Value = 2490.75
$ws.Cells.Item(113, 12).Value = 2687
$ws.Cells.Item(113, 13).Value = -320.75
$ws.Cells.Item(113, 14).Value = -7027
$ws.Cells.Item(122, 8).Value = 267282.38
$ws.Cells.Item(122, 9).Value = 553082.4
$ws.Cells.Item(122, 10).Value = 7464.227
$ws.Cells.Item(122, 11).Value = 1659247.2
$ws.Cells.Item(122, 12).Value = 22392.681
$ws.Cells.Item(122, 13).Value = -1656797.2
$ws.Cells.Item(122, 14).Value = -27292.681
$ws.Cells.Item(123, 8).Value = 53257.453
$ws.Cells.Item(123, 10).Value = 53257.453
$ws.Cells.Item(123, 12).Value = 53257.453
$ws.Cells.Item(123, 14).Value = -58157.453
$ws.Cells.Item(126, 8).Value = 7966.4287
$ws.Cells.Item(126, 9).Value = 2626.6667
$ws.Cells.Item(126, 11).Value = 7880.000100000001
$ws.Cells.Item(126, 13).Value = -5410.000100000001
$ws.Cells.Item(132, 8).Value = 9580
$ws.Cells.Item(132, 9).Value = 9300
$ws.Cells.Item(132, 11).Value = 27900
$ws.Cells.Item(132, 13).Value = -25370

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 4783.8213
$ws.Cells.Item(7, 9).Value = 3308.3572
$ws.Cells.Item(7, 11).Value = 3308.3572
$ws.Cells.Item(7, 13).Value = -3196.3572
$ws.Cells.Item(22, 8).Value = 1137.6923
$ws.Cells.Item(22, 9).Value = 956.4286
$ws.Cells.Item(22, 10).Value = 1349.1666
$ws.Cells.Item(22, 11).Value = 956.4286
$ws.Cells.Item(22, 12).Value = 1349.1666
$ws.Cells.Item(22, 13).Value = -661.4286
$ws.Cells.Item(22, 14).Value = -1939.1666
$ws.Cells.Item(27, 8).Value = 1137.6923
$ws.Cells.Item(27, 9).Value = 956.4286
$ws.Cells.Item(27, 10).Value = 1349.1666
$ws.Cells.Item(27, 11).Value = 956.4286
$ws.Cells.Item(27, 12).Value = 1349.1666
$ws.Cells.Item(27, 13).Value = -849.4286
$ws.Cells.Item(27, 14).Value = -1563.1666
$ws.Cells.Item(40, 8).Value = 13075691
$ws.Cells.Item(40, 9).Value = 4447287
$ws.Cells.Item(40, 10).Value = 37043480
$ws.Cells.Item(40, 11).Value = 4447287
$ws.Cells.Item(40, 12).Value = 37043480
$ws.Cells.Item(40, 13).Value = -4447151
$ws.Cells.Item(40, 14).Value = -37043752
$ws.Cells.Item(122, 8).Value = 16852.928
$ws.Cells.Item(122, 9).Value = 70000
$ws.Cells.Item(122, 11).Value = 210000
$ws.Cells.Item(122, 13).Value = -207550
$ws.Cells.Item(126, 8).Value = 4783.8213
$ws.Cells.Item(126, 9).Value = 3308.3572
$ws.Cells.Item(126, 11).Value = 9925.071599999999
$ws.Cells.Item(126, 13).Value = -7455.071599999999
$ws.Cells.Item(127, 8).Value = 250081420
$ws.Cells.Item(127, 10).Value = 108571.664
$ws.Cells.Item(127, 12).Value = 108571.664
$ws.Cells.Item(127, 14).Value = -118491.664
$ws.Cells.Item(132, 8).Value = 6214
$ws.Cells.Item(132, 9).Value = 5221
$ws.Cells.Item(132, 10).Value = 8200
$ws.Cells.Item(132, 11).Value = 15663
$ws.Cells.Item(132, 12).Value = 24600
$ws.Cells.Item(132, 13).Value = -13133
$ws.Cells.Item(132, 14).Value = -29660
$ws.Cells.Item(134, 8).Value = 81629
$ws.Cells.Item(134, 10).Value = 81629
$ws.Cells.Item(134, 12).Value = 81629
$ws.Cells.Item(134, 14).Value = -91769
$ws.Cells.Item(136, 8).Value = 5180.3125
$ws.Cells.Item(136, 9).Value = 2576.4
$ws.Cells.Item(136, 11).Value = 7729.200000000001
$ws.Cells.Item(136, 13).Value = -5179.200000000001

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 3477037.8
$ws.Cells.Item(81, 9).Value = 6945410.5
$ws.Cells.Item(81, 10).Value = 8665
$ws.Cells.Item(81, 11).Value = 13890821
$ws.Cells.Item(81, 12).Value = 17330
$ws.Cells.Item(81, 13).Value = -13889760
$ws.Cells.Item(81, 14).Value = -19452
$ws.Cells.Item(84, 8).Value = 3477037.8
$ws.Cells.Item(84, 9).Value = 6945410.5
$ws.Cells.Item(84, 10).Value = 8665
$ws.Cells.Item(84, 11).Value = 69454105
$ws.Cells.Item(84, 12).Value = 86650
$ws.Cells.Item(84, 13).Value = -69448801
$ws.Cells.Item(84, 14).Value = -97258
$ws.Cells.Item(122, 8).Value = 7199.9375
$ws.Cells.Item(122, 9).Value = 6463.125
$ws.Cells.Item(122, 11).Value = 19389.375
$ws.Cells.Item(122, 13).Value = -16939.375
$ws.Cells.Item(126, 8).Value = 6500.143
$ws.Cells.Item(126, 9).Value = 5464.6665
$ws.Cells.Item(126, 10).Value = 7276.75
$ws.Cells.Item(126, 11).Value = 16393.9995
$ws.Cells.Item(126, 12).Value = 21830.25
$ws.Cells.Item(126, 13).Value = -13923.9995
$ws.Cells.Item(126, 14).Value = -26770.25
$ws.Cells.Item(132, 8).Value = 55560776
$ws.Cells.Item(132, 9).Value = 2396
$ws.Cells.Item(132, 10).Value = 71434600
$ws.Cells.Item(132, 11).Value = 7188
$ws.Cells.Item(132, 12).Value = 214303800
$ws.Cells.Item(132, 13).Value = -4658
$ws.Cells.Item(132, 14).Value = -214308860
$ws.Cells.Item(136, 8).Value = 7297.97
$ws.Cells.Item(136, 9).Value = 3130.4888
$ws.Cells.Item(136, 10).Value = 10707.728
$ws.Cells.Item(136, 11).Value = 9391.466400000001
$ws.Cells.Item(136, 12).Value = 32123.184
$ws.Cells.Item(136, 13).Value = -6841.466400000001
$ws.Cells.Item(136, 14).Value = -37223.18399999999
